$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix config's bug: the second header cell was "ID" (duplicate/typo of "Id").
# Correct it to match cell A1 ("Id").
$ws.Range("B1").Value = "Id"

# Move the active selection to F10 (cursor position after the edit).
$ws.Range("F10").Select() | Out-Null
